$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the values for the model inputs in column C
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 0.15
$ws.Range("C4").Value = 0.4
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 3000

# Update the current selection to match the saved state
$ws.Range("H11").Select()
